$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.124.75'
$ws.Range("E2").Value = '  +0.45%  '
$ws.Range("D3").Value = '1.918.94'
$ws.Range("E3").Value = '  +2.55%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '319.69'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.08%  '
$ws.Range("E6").Value = '  +0.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5063'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.32%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4074'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.51%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08343'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.92%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.34'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.35%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.110'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.47%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.03'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.57%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.426'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.48%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.914.74'
$ws.Range("E14").Value = '  +2.56%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.247'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.005'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.26%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.55'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.67%  '
$ws.Range("E18").Value = '  +1.22%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06519'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.45%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.54'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.52%  '
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.947'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.43%  '
$ws.Range("D23").Value = '30.137.14'
$ws.Range("E23").Value = '  +0.53%  '
$ws.Range("E24").Value = '  +2.43%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.193'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.17%  '
$ws.Range("D26").Value = '2.135.42'
$ws.Range("E26").Value = '  +2.53%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.88'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.46%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '162.61'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.93%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.269'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.35%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '128.98'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.50%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.135'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +7.41%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1046'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.11%  '
$ws.Range("E33").Value = '  +0.89%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.788'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.26%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02458'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.66%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.329'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.47%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06440'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.61%  '
$ws.Range("E38").Value = '  +0.52%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6519'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.43%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.199'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.20%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.596'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.44'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.93%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.210'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.18%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.41'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.88%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.194'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +9.99%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6064'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.74%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.625'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.27%  '
$ws.Range("E48").Value = '  +0.79%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '122.31'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.29%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.146'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.52%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.98'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.07%  '
